$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Add the new "Est. Time" column (I) header and its value for the first
#    (My Requests) story block.
# ---------------------------------------------------------------------------
$ws.Range("A4").Copy() | Out-Null
$ws.Range("I2").PasteSpecial(-4122) | Out-Null
$ws.Range("I2").Value = "Est. Time"

$ws.Range("A4").Copy() | Out-Null
$ws.Range("I13").PasteSpecial(-4122) | Out-Null
$ws.Range("I13").Value = "4 Hours"

# ---------------------------------------------------------------------------
# 2. Rebuild the "Events Landing / Event details" block (rows 17-32).
#    Clear out everything from the old placeholder rows first so stray
#    cells/styles from the previous layout don't linger around.
# ---------------------------------------------------------------------------
$ws.Range("A17:H30").Clear() | Out-Null

# -- Row 17: section header -------------------------------------------------
$ws.Range("A4").Copy() | Out-Null
$ws.Range("A17").PasteSpecial(-4122) | Out-Null
$ws.Range("A17").Value = 2

$ws.Range("B4").Copy() | Out-Null
$ws.Range("B17").PasteSpecial(-4122) | Out-Null
$ws.Range("B17").Value = "Events Landing  & Event details"

# -- Rows 18-26: stories / tasks --------------------------------------------
$ws.Range("C18:C26").Clear() | Out-Null
$ws.Range("A4").Copy() | Out-Null
$ws.Range("C18:C26").PasteSpecial(-4122) | Out-Null

$ws.Range("D5").Copy() | Out-Null
$ws.Range("D18:D26").PasteSpecial(-4122) | Out-Null
$ws.Range("E18:E20").PasteSpecial(-4122) | Out-Null
$ws.Range("F18:F20").PasteSpecial(-4122) | Out-Null

$ws.Range("C18").Value = "US-01"
$ws.Range("D18").Value = 'As a user, I want to click on an "Events" link to see a list of upcoming events. (Events Landing)'
$ws.Range("E18").Value = "T-01"
$ws.Range("F18").Value = "Develop API calls to fetch event data from events_category-listing and events_listing endpoints."

$ws.Range("C19").Value = "US-02"
$ws.Range("D19").Value = "As a user, I want the Events Landing page to display event details like name, date, image, category, description, and location. (Events Landing)"
$ws.Range("E19").Value = "T-02"
$ws.Range("F19").Value = "Design and implement the Events Landing page with filtering and sorting functionalities."

$ws.Range("C20").Value = "US-03"
$ws.Range("D20").Value = " As a user, I want to see a maximum of 6 upcoming events with paging enabled for older events. (Events Landing)"
$ws.Range("E20").Value = "T-03"
$ws.Range("F20").Value = "Design and implement the Event Details page with registration button." + "`n" + "Implement logic to check for past events and disable registration."

$ws.Range("C21").Value = "US-04"
$ws.Range("D21").Value = " As a user, I want to filter events by date, category, or available tickets. (Events Landing)"

$ws.Range("C22").Value = "US-05"
$ws.Range("D22").Value = "As a user, I want events to be sorted by date in ascending order. (Events Landing)"

$ws.Range("C23").Value = "US-06"
$ws.Range("D23").Value = "As a user, I want to click on an event name to see its detailed information on a separate page. (Event Details)"

$ws.Range("C24").Value = "US-07"
$ws.Range("D24").Value = "As a user, I want the Event Details page to display additional information like end date and event tickets. (Event Details)"

$ws.Range("C25").Value = "US-08"
$ws.Range("D25").Value = 'As a user, I want a "Register Now" button to register for an event, redirecting me to the registration page. (Event Details)'

$ws.Range("C26").Value = "US-09"
$ws.Range("D26").Value = " As a user, I shouldn't be able to register for events that have already passed. (Event Details)"

# -- Row heights for the new content rows -----------------------------------
$ws.Rows(18).RowHeight = 30
$ws.Rows(19).RowHeight = 45
$ws.Rows(20).RowHeight = 60
$ws.Rows(21).RowHeight = 30
$ws.Rows(22).RowHeight = 30
$ws.Rows(23).RowHeight = 30
$ws.Rows(24).RowHeight = 30
$ws.Rows(25).RowHeight = 30
$ws.Rows(26).RowHeight = 30

# -- Row 27: empty separator row (keeps style of column C) ------------------
$ws.Range("A4").Copy() | Out-Null
$ws.Range("C27").PasteSpecial(-4122) | Out-Null

# -- Row 28: Workspace Dashboard ---------------------------------------------
$ws.Range("A4").Copy() | Out-Null
$ws.Range("A28").PasteSpecial(-4122) | Out-Null
$ws.Range("A28").Value = 3

$ws.Range("B4").Copy() | Out-Null
$ws.Range("B28").PasteSpecial(-4122) | Out-Null
$ws.Range("B28").Value = "Workspace Dashboard"

# -- Row 32: Registration Request details page -------------------------------
$ws.Range("A4").Copy() | Out-Null
$ws.Range("A32").PasteSpecial(-4122) | Out-Null
$ws.Range("A32").Value = 3

$ws.Range("B4").Copy() | Out-Null
$ws.Range("B32").PasteSpecial(-4122) | Out-Null
$ws.Range("B32").Value = "Registration Request details page"

# ---------------------------------------------------------------------------
# 3. Column I width
# ---------------------------------------------------------------------------
$ws.Columns("I").ColumnWidth = 22.86

# ---------------------------------------------------------------------------
# 4. Selection / view state: select row 27 (matches the entry the author was
#    looking at when they saved the file).
# ---------------------------------------------------------------------------
$ws.Rows(27).Select() | Out-Null

$wb.Application.CutCopyMode = $false
